$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 3402
$ws.Range("K3").Value = 3376
$ws.Range("F4").Value = 1909
$ws.Range("H4").Value = 1731
$ws.Range("J4").Value = 1821
$ws.Range("K4").Value = 707
$ws.Range("K5").Value = 222
$ws.Range("K6").Value = 3968
$ws.Range("F7").Value = 24102
$ws.Range("H7").Value = 26044
$ws.Range("J7").Value = 29292
$ws.Range("K7").Value = 11675

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K3").Value = 33
$ws.Range("K6").Value = 79
$ws.Range("K7").Value = 155

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 226
$ws.Range("K3").Value = 233
$ws.Range("K6").Value = 255
$ws.Range("K7").Value = 778

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K3").Value = 81
$ws.Range("K6").Value = 61
$ws.Range("K7").Value = 251

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 130
$ws.Range("K3").Value = 175
$ws.Range("K6").Value = 132
$ws.Range("K7").Value = 469

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 61
$ws.Range("K3").Value = 68
$ws.Range("K7").Value = 191

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K3").Value = 140
$ws.Range("K7").Value = 404

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K6").Value = 110
$ws.Range("K7").Value = 282

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K2").Value = 59
$ws.Range("K7").Value = 203

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("K4").Value = 1
$ws.Range("K7").Value = 37

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K4").Value = 41
$ws.Range("K7").Value = 331
$ws.Range("K8").Value = 778
$ws.Range("K11").Value = 243
$ws.Range("K14").Value = 60
$ws.Range("K19").Value = 359
$ws.Range("K20").Value = 269
$ws.Range("K27").Value = 119
$ws.Range("F29").Value = 1655
$ws.Range("K29").Value = 612
$ws.Range("K30").Value = 37
$ws.Range("K33").Value = 469
$ws.Range("K36").Value = 136
$ws.Range("K37").Value = 404
$ws.Range("K40").Value = 30
$ws.Range("K41").Value = 102
$ws.Range("K42").Value = 417
$ws.Range("K43").Value = 106
$ws.Range("K46").Value = 24
$ws.Range("K47").Value = 64
$ws.Range("K51").Value = 140
$ws.Range("K52").Value = 317
$ws.Range("K53").Value = 155
$ws.Range("K54").Value = 230
$ws.Range("K55").Value = 123
$ws.Range("K60").Value = 72
$ws.Range("H63").Value = 283
$ws.Range("J63").Value = 103
$ws.Range("K63").Value = 38
$ws.Range("K65").Value = 282
$ws.Range("K67").Value = 454
$ws.Range("K71").Value = 34
$ws.Range("K73").Value = 104
$ws.Range("K76").Value = 177
$ws.Range("K77").Value = 83
$ws.Range("K78").Value = 147
$ws.Range("K79").Value = 302
$ws.Range("K83").Value = 251
$ws.Range("K85").Value = 544
$ws.Range("K86").Value = 79
$ws.Range("K89").Value = 156
$ws.Range("K91").Value = 121
$ws.Range("K92").Value = 45
$ws.Range("K94").Value = 144
$ws.Range("K95").Value = 191
$ws.Range("K96").Value = 144
$ws.Range("K99").Value = 203
$ws.Range("F101").Value = 24102
$ws.Range("H101").Value = 26044
$ws.Range("J101").Value = 29292
$ws.Range("K101").Value = 11675

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 148
$ws.Range("K6").Value = 133
$ws.Range("K7").Value = 454

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K3").Value = 69
$ws.Range("K4").Value = 12
$ws.Range("K6").Value = 106
$ws.Range("K7").Value = 230

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 170
$ws.Range("F4").Value = 92
$ws.Range("K4").Value = 35
$ws.Range("K5").Value = 15
$ws.Range("K6").Value = 186
$ws.Range("F7").Value = 1655
$ws.Range("K7").Value = 612

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K3").Value = 96
$ws.Range("K6").Value = 114
$ws.Range("K7").Value = 359

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K6").Value = 103
$ws.Range("K7").Value = 177

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("K2").Value = 26
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K2").Value = 34
$ws.Range("K4").Value = 8
$ws.Range("K7").Value = 102

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K6").Value = 155
$ws.Range("K7").Value = 417

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K3").Value = 32
$ws.Range("K7").Value = 147

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 123

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("K2").Value = 10
$ws.Range("K7").Value = 24

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K3").Value = 39
$ws.Range("K6").Value = 28

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K3").Value = 25
$ws.Range("K7").Value = 144

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K2").Value = 32
$ws.Range("K7").Value = 121

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 104
$ws.Range("K4").Value = 18
$ws.Range("K7").Value = 302

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K3").Value = 76
$ws.Range("K6").Value = 87
$ws.Range("K7").Value = 269

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 136

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 121
$ws.Range("K6").Value = 85
$ws.Range("K7").Value = 331

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K6").Value = 61
$ws.Range("K7").Value = 144

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K3").Value = 22
$ws.Range("K7").Value = 64

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 73
$ws.Range("K7").Value = 243

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K6").Value = 43
$ws.Range("K7").Value = 104

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("K2").Value = 14
$ws.Range("K7").Value = 45

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K2").Value = 39
$ws.Range("K3").Value = 51
$ws.Range("K7").Value = 156

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K3").Value = 29
$ws.Range("K7").Value = 119

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K4").Value = 30
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 79

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K2").Value = 37
$ws.Range("K3").Value = 38
$ws.Range("K4").Value = 13
$ws.Range("K7").Value = 140

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K3").Value = 24
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 72

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K2").Value = 22
$ws.Range("K6").Value = 45
$ws.Range("K7").Value = 106

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 193
$ws.Range("K7").Value = 544

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("K2").Value = 17
$ws.Range("K7").Value = 34

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 83

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 30

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K3").Value = 82
$ws.Range("K7").Value = 317

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 41
